$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns D..AQ, applied identically to rows 2 and 3
$values = @{
    "D"  = 0.0803
    "E"  = -0.06509999999999999
    "F"  = 0.115
    "G"  = 0.1283299180327869
    "H"  = 0.1283299180327869
    "I"  = 0.1196666422716628
    "J"  = 0.09993423432769508
    "K"  = 1411.5
    "L"  = 0.1291258050351288
    "M"  = 978.2
    "N"  = 0.04165318253819557
    "O"  = 0.6930216082182076
    "P"  = 978.2
    "Q"  = 0.04165318253819557
    "R"  = 0.6930216082182076
    "U"  = 1380.9
    "V"  = 0.05880073580759994
    "W"  = 0.1148513401357223
    "X"  = 0.05052585406158823
    "Y"  = 0.06432548607413409
    "Z"  = 0.9811158181948735
    "AA" = 0.09804705807809476
    "AB" = 0.04382671695860546
    "AC" = 0.0542203411194893
    "AD" = 6516
    "AE" = 0
    "AF" = 6516
    "AG" = 5135.1
    "AH" = 0.2171971040386128
    "AI" = 0.3105903886669781
    "AJ" = 0.1794266147207324
    "AK" = 0.262014950123734
    "AL" = 83.40000000000001
    "AM" = 83.40000000000001
    "AN" = 4.766293614219882
    "AO" = 15.68465227817746
    "AP" = 3.756199253895107
    "AQ" = 15.68465227817746
}

foreach ($row in 2..3) {
    foreach ($col in $values.Keys) {
        $addr = "$col$row"
        $ws.Range($addr).Value = $values[$col]
    }
}
